$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: fix header text that shifts up after crystal_nova_x.lua string removal ---
# (G1/H1 already hold "技能键值"/"键值描述" - no value change needed, just keep as-is)

# --- Add new header "预载资源" / "Precache" in column AJ (row1/row2) ---
$ws.Range("AJ1").Value = "预载资源"
$ws.Range("AJ2").Value = "Precache"

# Copy the header style (row1 uses style like G1/H1, row2 uses style like other row2 cells)
$ws.Range("G1").Copy()
$ws.Range("AJ1").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy()
$ws.Range("AJ2").PasteSpecial(-4122) | Out-Null

# --- Row 3: E3 becomes a formula instead of a literal string ---
$ws.Range("E3").Formula = "=""examples/abilities/""&A3&"".lua"""

# --- Row 4: new ability example "counter_helix_x" (axe counter helix) ---
$ws.Range("A4").Value = "counter_helix_x"
$ws.Range("B4").Value = "反击螺旋"
$ws.Range("C4").Value = "受到一定次数攻击后，斧王就会做出螺旋反击，对附近所有敌方单位造成纯粹伤害。"
$ws.Range("D4").Value = "ability_lua"
$ws.Range("E4").Formula = "=""examples/abilities/""&A4&"".lua"""

$ws.Range("G4").Value = "damage 75 110 145 180"
$ws.Range("H4").Value = "伤害"
$ws.Range("I4").Value = "radius 275"
$ws.Range("J4").Value = "范围"
$ws.Range("K4").Value = "hit_count 7 6 5 4"
$ws.Range("L4").Value = "所需攻击次数"

$ws.Range("X4").Value = 4

$ws.Range("AA4").Value = "axe_counter_helix"

$jsonText = "{" + [char]10 + [char]9 + '"soundfile"' + [char]9 + '"soundevents/game_sounds_heroes/game_sounds_axe.vsndevts"' + [char]10 + [char]9 + '"particle"' + [char]9 + '"particles/units/heroes/hero_axe/axe_counterhelix_ad.vpcf"' + [char]10 + "}"
$ws.Range("AJ4").Value = $jsonText

# --- Copy row 3 formatting down into row 4 so styles (borders/fonts) match ---
$ws.Range("A3:L3").Copy()
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("X3").Copy()
$ws.Range("X4").PasteSpecial(-4122) | Out-Null
$ws.Range("AA3").Copy()
$ws.Range("AA4").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(4).RowHeight = 14.25

# AJ4 needs wrap text (new style entry) for the multi-line precache JSON block
$ws.Range("AJ4").WrapText = $true

# --- Selection / view state ---
$ws.Range("AA4").Select()

$excel.CutCopyMode = $false
